# Generate Report for Handoff
# Rows in "Overview", "zh-cn" and "de-de" that were "Ready for handoff"
# (sharing the 2016-08-28 12:21:xx handoff timestamps) now get a fresh
# xliff handoff generated: their "Priority" column is flagged "ht" and
# the handoff timestamps are bumped forward.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 13)

# --- Overview sheet: bump "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-28 12:22:09"
}

# --- zh-cn sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-28 12:22:00"
}

# --- de-de sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-28 12:22:09"
}
